$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the hourly crypto price/volume snapshot (Price and Volume(1h)
# columns) for rows 2-51. Values that look like plain decimal numbers are
# written with a leading apostrophe so Excel keeps them as text (matching
# the original formatting) instead of normalizing them into floating point
# numbers and losing trailing zeros / precision.
$ws.Range("D2").Value = "26.937.16"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.847.69"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'309.70"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4770"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("D8").Value = "'0.3686"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").Value = "'0.07197"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "'0.9266"
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("D11").Value = "'19.63"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'0.07615"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").Value = "1.855.44"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "'5.307"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "'6.403"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'88.54"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'0.000008647"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "26.972.81"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "'14.56"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").Value = "'5.027"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'1.918"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'152.38"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").Value = "'2.010"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").Value = "'114.35"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "'4.915"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").Value = "'0.08868"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "'3.277"
$ws.Range("E31").Value = "  +4.35%  "
$ws.Range("D32").Value = "'0.7485"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").Value = "'1.166"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("D34").Value = "'2.787"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'4.495"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").Value = "'0.05264"
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "'2.973"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "'0.5225"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").Value = "'6.903"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "'0.1514"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D43").Value = "'8.214"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").Value = "'10.54"
$ws.Range("E44").Value = "  +5.00%  "
$ws.Range("D45").Value = "'0.4703"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'102.06"
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("D48").Value = "'1.606"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "'65.46"
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("D50").Value = "'0.06038"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "'0.8853"
$ws.Range("E51").Value = "  +4.44%  "
